# "Generate Report for Handback"
#
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns on the per-language report sheets (zh-cn,
# de-de) once a handback has completed, flips the Overview status text
# from "In Translation" to "Handed back: in sync with en-US", and widens
# the columns that now hold longer text so it is readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (shared by Overview!E2:F3 and the "Status" column on both language
#    sheets, since they all point at the same shared string).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 2. Per-language sheets: record the handed-back target files and the
#    handback timestamps.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$md1 = "4d6a1377-b97f-4ed7-9693-2b138a410528.md"
$md2 = "4fa7bead-87d8-4a89-87be-4465f701b10d.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b0d2ee36ba9931a38a42771f4e6502ad8e2fb4b/e2e/4d6a1377-b97f-4ed7-9693-2b138a410528.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b0d2ee36ba9931a38a42771f4e6502ad8e2fb4b/e2e/4fa7bead-87d8-4a89-87be-4465f701b10d.md"

# zh-cn
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("I2").Value = $md1
$zhcn.Range("J2").Value = "4d6a1377-b97f-4ed7-9693-2b138a410528.e9ad60a49cec07b244601e641a00c3ec1c5fa31d.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 08:27:17"
$zhcn.Range("I3").Value = $md2
$zhcn.Range("J3").Value = "4fa7bead-87d8-4a89-87be-4465f701b10d.09fb933d779e22013c3fb037113400d0b76eaf77.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-30 08:27:17"

# de-de
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("I2").Value = $md1
$dede.Range("J2").Value = "4d6a1377-b97f-4ed7-9693-2b138a410528.e9ad60a49cec07b244601e641a00c3ec1c5fa31d.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 08:27:25"
$dede.Range("I3").Value = $md2
$dede.Range("J3").Value = "4fa7bead-87d8-4a89-87be-4465f701b10d.09fb933d779e22013c3fb037113400d0b76eaf77.de-de.xlf"
$dede.Range("K3").Value = "2016-08-30 08:27:25"

# ---------------------------------------------------------------------
# 3. New hyperlinks on the "Latest Target File" cells, re-using the same
#    targets as the existing "Source File Name" hyperlinks. Rebuild all
#    hyperlinks on each sheet so relationship ids come out in the same
#    (document) order as Excel itself would assign them.
# ---------------------------------------------------------------------
foreach ($ws in @($zhcn, $dede)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $url1, "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $url2, "", "", $md2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, "", "", $md2)
}

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer text. Excel quantizes
#    ColumnWidth to whole pixels, so request the character width whose
#    rounded pixel width is closest to the intended value.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.2   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.2   # F: de-de status

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth = 29.2      # C: Status
    $ws.Columns.Item(9).ColumnWidth = 39.2      # I: Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.2     # J: Latest Handback File
}
